# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document.") as three *separate* runs:
#   " (" / "Changed main" / ")"
# (matching the target OOXML diff, which shows 4 distinct <w:r> elements
#  inside the existing first paragraph).
#
# Plain Range.InsertAfter() on this host coalesces new text into the
# neighbouring run when the formatting is identical, so instead we stage the
# new text in a scratch paragraph at the end of the story, Copy() each piece
# and Paste() it at the end of paragraph 1 - paste always lands as its own
# run and never merges with what is already there. The scratch paragraph is
# then removed again so the rest of the document is unaffected.

$d = $word.ActiveDocument

$targetParagraph = $d.Paragraphs(1)

# --- 1. Build a scratch paragraph at the very end of the document holding
#        the text we need to copy from. -------------------------------------
$docEnd = $d.Content.End
$tailRange = $d.Range($docEnd - 1, $docEnd - 1)
$tailRange.InsertParagraphAfter()

$scratchParaIndex = $d.Paragraphs.Count
$scratchParagraph = $d.Paragraphs($scratchParaIndex)
$scratchRange = $scratchParagraph.Range
$scratchRange.Collapse(1)
$scratchRange.InsertAfter(" (Changed main)")

# Piece boundaries (character offsets from the scratch paragraph start)
# within the scratch text " (Changed main)":
#   " ("           -> [0, 2)
#   "Changed main" -> [2, 14)
#   ")"            -> [14, 15)
$pieceBounds = @(
    @(0, 2),
    @(2, 14),
    @(14, 15)
)

# --- 2. Copy each piece and paste it onto the end of paragraph 1. ----------
# The scratch paragraph is re-resolved by index on every iteration because
# its character offsets shift forward each time text is pasted earlier in
# the document (into paragraph 1).
foreach ($bounds in $pieceBounds) {
    $liveScratch = $d.Paragraphs($scratchParaIndex)
    $liveStart = $liveScratch.Range.Start

    $pieceStart = $liveStart + $bounds[0]
    $pieceEnd = $liveStart + $bounds[1]
    $piece = $d.Range($pieceStart, $pieceEnd)
    $piece.Copy()

    $insertAt = $targetParagraph.Range.End
    $dest = $d.Range($insertAt - 1, $insertAt - 1)
    $dest.Paste()
}

# --- 3. Remove the scratch paragraph (text + its paragraph mark). ----------
$scratchParagraph = $d.Paragraphs($scratchParaIndex)
$scratchFull = $d.Range($scratchParagraph.Range.Start, $scratchParagraph.Range.End)
$scratchFull.Delete()

Write-Output "Paragraph 1 is now: $($targetParagraph.Range.Text)"
